$p = $ppt.ActivePresentation

# Both text edits described by the commit ("updated System Design graphic")
# live on the slide with the "Mary Lou ... System Design" title.
$s2 = $p.Slides.Item(2)

# --- Locate the shapes we need to touch -------------------------------------
$attitudeBox = $null
$titleShape = $null

for ($i = 1; $i -le $s2.Shapes.Count; $i++) {
    $shp = $s2.Shapes.Item($i)

    if ($shp.Type -eq 6) {
        # a group shape - look one level down for the "Attitude" label
        $items = $shp.GroupItems
        for ($j = 1; $j -le $items.Count; $j++) {
            $sub = $items.Item($j)
            if ($sub.HasTextFrame -and $sub.TextFrame.HasText -and $sub.TextFrame.TextRange.Text -eq "Attitude") {
                $attitudeBox = $sub
            }
        }
    } elseif ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
        if ($shp.TextFrame.TextRange.Text -eq "Mary Lou Overall System Design") {
            $titleShape = $shp
        }
    }
}

# --- 1. "Attitude" -> "Pitch" ------------------------------------------------
if ($attitudeBox -ne $null) {
    $attitudeBox.TextFrame.TextRange.Text = "Pitch"
}

# --- 2. Title: "Mary Lou Overall System Design" -> "Mary Lou System Design" -
# Final text is split across 4 runs: "Mary ", "Lou ", "System ", "Design"
if ($titleShape -ne $null) {
    $tr = $titleShape.TextFrame.TextRange

    # Split "Mary Lou " into "Mary " + "Lou "
    $tr.Characters(6, 4).Text = "Lou "

    # Collapse "Overall System Design" down to "System Design"
    $tr.Characters(10, 22).Text = "System Design"

    # Split "System Design" into "System " + "Design"
    $tr.Characters(17, 6).Text = "Design"
}
